# South_unemployment.xlsx update
# - Revise D217:D276 values (2017 M01 .. 2021 M12) per BLS March 2022 revision
# - Append two new rows for 2022 M01 and M02
# - Add a revision comment to each revised cell in D217:D276
# - Update the "Years:" label from "2000 to 2021" to "2000 to 2022"
# - Update the footer generation timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Revised monthly values for rows 217-276 (2017 M01 through 2021 M12)
# ---------------------------------------------------------------------------
$revisedValues = @{
    217 = 2788626; 218 = 2736500; 219 = 2681931; 220 = 2630922; 221 = 2583849;
    222 = 2545382; 223 = 2518480; 224 = 2496887; 225 = 2476232; 226 = 2454733;
    227 = 2433804; 228 = 2410481; 229 = 2386621; 230 = 2364041; 231 = 2340808;
    232 = 2314424; 233 = 2288816; 234 = 2264886; 235 = 2242477; 236 = 2229336;
    237 = 2229566; 238 = 2238219; 239 = 2246296; 240 = 2244592; 241 = 2227144;
    242 = 2192308; 243 = 2151434; 244 = 2114920; 245 = 2091284; 246 = 2085799;
    247 = 2095797; 248 = 2108602; 249 = 2107588; 250 = 2097685; 251 = 2087160;
    252 = 2084382; 253 = 2092659; 254 = 2111509; 255 = 2579147; 256 = 7454191;
    257 = 6758450; 258 = 5690870; 259 = 5425345; 260 = 4481515; 261 = 4256588;
    262 = 3793598; 263 = 3623935; 264 = 3543818; 265 = 3384250; 266 = 3292703;
    267 = 3196300; 268 = 3147187; 269 = 3043215; 270 = 3014585; 271 = 2877264;
    272 = 2769653; 273 = 2614951; 274 = 2523062; 275 = 2455032; 276 = 2429637
}

$commentText = "*  Data were subject to revision on March 2, 2022.`n"

for ($r = 217; $r -le 276; $r++) {
    $ws.Cells.Item($r, 4).Value = $revisedValues[$r]
    $ws.Cells.Item($r, 4).AddComment($commentText) | Out-Null
}

# ---------------------------------------------------------------------------
# 2. Append two new rows of data for 2022 (M01, M02), copying the format of
#    the last existing data row so styles line up with the rest of the table
# ---------------------------------------------------------------------------
$ws.Range("A276:D276").Copy() | Out-Null
$ws.Range("A277:D278").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Cells.Item(277, 1).Value = "LASRD930000000000004"
$ws.Cells.Item(277, 2).Value = 2022
$ws.Cells.Item(277, 3).Value = "M01"
$ws.Cells.Item(277, 4).Value = 2378613

$ws.Cells.Item(278, 1).Value = "LASRD930000000000004"
$ws.Cells.Item(278, 2).Value = 2022
$ws.Cells.Item(278, 3).Value = "M02"
$ws.Cells.Item(278, 4).Value = 2320399

# ---------------------------------------------------------------------------
# 3. Update the "Years:" label
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = "2000 to 2022"

# ---------------------------------------------------------------------------
# 4. Update the footer generation timestamp
# ---------------------------------------------------------------------------
$ws.PageSetup.LeftFooter = "Source: Bureau of Labor Statistics"
$ws.PageSetup.RightFooter = "Generated on: March 28, 2022 (06:25:15 PM)"
